# Scheduled-runner update: refresh market-board snapshot values (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 14298780
$ws.Range("I69").Value = 50005980
$ws.Range("K69").Value = 150017940
$ws.Range("M69").Value = -150017066

$ws.Range("H72").Value = 14298780
$ws.Range("I72").Value = 50005980
$ws.Range("K72").Value = 450053820
$ws.Range("M72").Value = -450049452

$ws.Range("H86").Value = 121216170
$ws.Range("I86").Value = 250003730
$ws.Range("J86").Value = 47623280
$ws.Range("K86").Value = 250003730
$ws.Range("L86").Value = 47623280
$ws.Range("M86").Value = -250002607
$ws.Range("N86").Value = -47625526

$ws.Range("H89").Value = 121216170
$ws.Range("I89").Value = 250003730
$ws.Range("J89").Value = 47623280
$ws.Range("K89").Value = 1250018650
$ws.Range("L89").Value = 238116400
$ws.Range("M89").Value = -1250013034
$ws.Range("N89").Value = -238127632

$ws.Range("H98").Value = 2992
$ws.Range("J98").Value = 9995
$ws.Range("L98").Value = 9995
$ws.Range("N98").Value = -12991

$ws.Range("H122").Value = 2992
$ws.Range("J122").Value = 9995
$ws.Range("L122").Value = 29985
$ws.Range("N122").Value = -34885

$ws.Range("H135").Value = 304.09525
$ws.Range("I135").Value = 304.09525
$ws.Range("K135").Value = 2736.85725
$ws.Range("M135").Value = -201.85725

$ws.Range("H138").Value = 5187.02
$ws.Range("I138").Value = 6555.0713
$ws.Range("J138").Value = 3445.8635
$ws.Range("K138").Value = 19665.2139
$ws.Range("L138").Value = 10337.5905
$ws.Range("M138").Value = -14525.2139
$ws.Range("N138").Value = -20617.5905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 99499.5
$ws.Range("J123").Value = 99499.5
$ws.Range("L123").Value = 99499.5
$ws.Range("N123").Value = -109299.5

$ws.Range("H138").Value = 69998.336
$ws.Range("J138").Value = 69998.336
$ws.Range("L138").Value = 69998.336
$ws.Range("N138").Value = -80278.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2000.6666
$ws.Range("I36").Value = 2000.6666
$ws.Range("K36").Value = 2000.6666
$ws.Range("M36").Value = -1466.6666

$ws.Range("H94").Value = 1096
$ws.Range("I94").Value = 1155
$ws.Range("J94").Value = 299.5
$ws.Range("K94").Value = 1155
$ws.Range("L94").Value = 299.5
$ws.Range("M94").Value = -704
$ws.Range("N94").Value = -1201.5

$ws.Range("H107").Value = 28471.072
$ws.Range("I107").Value = 30276.54
$ws.Range("K107").Value = 30276.54
$ws.Range("M107").Value = -28356.54

$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 9999.666999999999
$ws.Range("I8").Value = 10000
$ws.Range("J8").Value = 9999.5
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 9999.5
$ws.Range("M8").Value = -9860
$ws.Range("N8").Value = -10279.5

$ws.Range("H58").Value = 2517.25
$ws.Range("I58").Value = 1337.3
$ws.Range("K58").Value = 1337.3
$ws.Range("M58").Value = -1134.3

$ws.Range("H105").Value = 3402.4
$ws.Range("J105").Value = 4750
$ws.Range("L105").Value = 4750
$ws.Range("N105").Value = -8244

$ws.Range("H132").Value = 1872.9667
$ws.Range("I132").Value = 1771.7
$ws.Range("K132").Value = 5315.1
$ws.Range("M132").Value = -2785.1

$ws.Range("H136").Value = 2517.25
$ws.Range("I136").Value = 1337.3
$ws.Range("K136").Value = 4011.9
$ws.Range("M136").Value = -1461.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1174.3235
$ws.Range("I2").Value = 486.6154
$ws.Range("J2").Value = 1600.0476
$ws.Range("K2").Value = 2919.6924
$ws.Range("L2").Value = 9600.285600000001
$ws.Range("M2").Value = -2806.6924
$ws.Range("N2").Value = -9826.285600000001

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws.Range("H75").Value = 17860444
$ws.Range("J75").Value = 23813518
$ws.Range("L75").Value = 71440554
$ws.Range("N75").Value = -71442550

$ws.Range("H78").Value = 17860444
$ws.Range("J78").Value = 23813518
$ws.Range("L78").Value = 214321662
$ws.Range("N78").Value = -214331646

$ws.Range("H103").Value = 716683.5
$ws.Range("I103").Value = 1666882.1
$ws.Range("J103").Value = 4034.5
$ws.Range("K103").Value = 5000646.300000001
$ws.Range("L103").Value = 12103.5
$ws.Range("M103").Value = -4999767.300000001
$ws.Range("N103").Value = -13861.5

$ws.Range("H114").Value = 4767.9287
$ws.Range("I114").Value = 1175.2
$ws.Range("J114").Value = 13749.75
$ws.Range("K114").Value = 3525.6
$ws.Range("L114").Value = 41249.25
$ws.Range("M114").Value = -271.6000000000004
$ws.Range("N114").Value = -47757.25

$ws.Range("H122").Value = 6668645
$ws.Range("J122").Value = 5296
$ws.Range("L122").Value = 47664
$ws.Range("N122").Value = -52564

$ws.Range("H131").Value = 5052948
$ws.Range("I131").Value = 8266025
$ws.Range("J131").Value = 3827.1428
$ws.Range("K131").Value = 24798075
$ws.Range("L131").Value = 11481.4284
$ws.Range("M131").Value = -24793035
$ws.Range("N131").Value = -21561.4284

$ws.Range("H133").Value = 4336.222
$ws.Range("I133").Value = 2238.7144
$ws.Range("K133").Value = 6716.1432
$ws.Range("M133").Value = -1656.1432

$ws.Range("H134").Value = 3510.0667
$ws.Range("I134").Value = 1093
$ws.Range("J134").Value = 5625
$ws.Range("K134").Value = 3279
$ws.Range("L134").Value = 16875
$ws.Range("M134").Value = 1791
$ws.Range("N134").Value = -27015

$ws.Range("H136").Value = 8586.817999999999
$ws.Range("I136").Value = 3909.5833
$ws.Range("K136").Value = 11728.7499
$ws.Range("M136").Value = -6628.749899999999

$ws.Range("H137").Value = 2227.3572
$ws.Range("J137").Value = 2678.8
$ws.Range("L137").Value = 8036.400000000001
$ws.Range("N137").Value = -18236.4

$ws.Range("H138").Value = 3662.85
$ws.Range("I138").Value = 3733.0588
$ws.Range("K138").Value = 11199.1764
$ws.Range("M138").Value = -6059.1764

$ws.Range("H139").Value = 4312945
$ws.Range("I139").Value = 5001816
$ws.Range("K139").Value = 15005448
$ws.Range("M139").Value = -15000308

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 500
$ws.Range("M14").Value = -332
$ws.Range("N14").Value = -836

$ws.Range("H97").Value = 167560.42
$ws.Range("I97").Value = 250204.75
$ws.Range("J97").Value = 126238.25
$ws.Range("K97").Value = 250204.75
$ws.Range("L97").Value = 126238.25
$ws.Range("M97").Value = -249708.75
$ws.Range("N97").Value = -127230.25

$ws.Range("H113").Value = 8506.875
$ws.Range("J113").Value = 8015.75
$ws.Range("L113").Value = 8015.75
$ws.Range("N113").Value = -12355.75

$ws.Range("H126").Value = 1500.4
$ws.Range("I126").Value = 1500.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4501.200000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2031.200000000001
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 2917.625
$ws.Range("I132").Value = 2501.3
$ws.Range("K132").Value = 7503.900000000001
$ws.Range("M132").Value = -4973.900000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3949.8845
$ws.Range("I22").Value = 1231.1428
$ws.Range("J22").Value = 4951.5264
$ws.Range("K22").Value = 1231.1428
$ws.Range("L22").Value = 4951.5264
$ws.Range("M22").Value = -936.1428000000001
$ws.Range("N22").Value = -5541.5264

$ws.Range("H23").Value = 1672600
$ws.Range("I23").Value = 2503900
$ws.Range("K23").Value = 2503900
$ws.Range("M23").Value = -2503670

$ws.Range("H27").Value = 3949.8845
$ws.Range("I27").Value = 1231.1428
$ws.Range("J27").Value = 4951.5264
$ws.Range("K27").Value = 1231.1428
$ws.Range("L27").Value = 4951.5264
$ws.Range("M27").Value = -1124.1428
$ws.Range("N27").Value = -5165.5264

$ws.Range("H46").Value = 8157.773
$ws.Range("I46").Value = 23697.6
$ws.Range("K46").Value = 23697.6
$ws.Range("M46").Value = -23509.6

$ws.Range("H68").Value = 9500
$ws.Range("I68").Value = 9500
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 9500
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -8751
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 9500
$ws.Range("I71").Value = 9500
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 47500
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -43756
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3099.8
$ws.Range("I62").Value = 2999.6667
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 2999.6667
$ws.Range("L62").Value = 3250
$ws.Range("M62").Value = -2375.6667
$ws.Range("N62").Value = -4498

$ws.Range("H65").Value = 3099.8
$ws.Range("I65").Value = 2999.6667
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 14998.3335
$ws.Range("L65").Value = 16250
$ws.Range("M65").Value = -11878.3335
$ws.Range("N65").Value = -22490

$ws.Range("H126").Value = 25003658
$ws.Range("I126").Value = 41669932
$ws.Range("K126").Value = 125009796
$ws.Range("M126").Value = -125007326
